# NEXT Group edit to natural gas heat rate
# Update the Electricity Dispatch Logit Exponent value on the EDLE sheet
# and restore the "About" sheet as the active/selected tab.

$wb = $excel.ActiveWorkbook

$wsEdle = $wb.Worksheets.Item("EDLE")
$wsAbout = $wb.Worksheets.Item("About")

# Core data edit: the logit exponent value.
$wsEdle.Range("B2").Value = -3

# Move the EDLE sheet's cursor position before switching away from it,
# then make "About" the active sheet/tab (matches saved view state).
$wsEdle.Range("A29").Select()
$wsAbout.Activate()
